$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Host" / "endoparasitic" translation row (row 4), shifting
# all subsequent rows up by one.
$ws.Rows(4).Delete()

# Restore the (now out-of-range) active cell selection left behind by the
# row deletion in the original edit.
$ws.Range("B18").Select()
